$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.917049579578596
$ws.Range("D2").Value = 0.3690609465635661

$ws.Range("C3").Value = 1.138316513323244
$ws.Range("D3").Value = 0.2672331538984094

$ws.Range("C4").Value = 1.202863719698396
$ws.Range("D4").Value = 0.2418149928164184

$ws.Range("C5").Value = 3.00029286156353
$ws.Range("D5").Value = 0.006590145372918688

$ws.Range("C6").Value = 0.4168925843958579
$ws.Range("D6").Value = 0.6807967533298473

$ws.Range("C7").Value = 0.3585080726329952
$ws.Range("D7").Value = 0.7233804878264385

$ws.Range("C8").Value = 2.232068840476776
$ws.Range("D8").Value = 0.03611291339940559

$ws.Range("C9").Value = -0.08885205859202046
$ws.Range("D9").Value = 0.9300033801534284

$ws.Range("C10").Value = 1.951234482192545
$ws.Range("D10").Value = 0.06387801263430659
$ws.Range("G10").Value = "No"

$ws.Range("C11").Value = 2.20186174640901
$ws.Range("D11").Value = 0.03845769984566982
$ws.Range("G11").Value = "Sí"
